# "Generate Report for Handback"
#
# The de-de / zh-cn handback files are now in sync with en-US, so the
# localization-status report is regenerated to reflect that:
#   - Status changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + each language sheet)
#   - each language sheet's "Latest Handback DateTime" is refreshed to the
#     handback run that just completed
#   - the stale "handback file is not the latest" Error Detail is cleared
#     now that the handback is current

$wb = $excel.ActiveWorkbook

$statusNew = "Handed back: in sync with en-US"

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F2").Value = $statusNew
$overview.Range("E1:F1").ColumnWidth = 29.1

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusNew
$zhcn.Range("K2").Value = "2016-08-26 06:47:41"
$zhcn.Range("P2").Value = ""
$zhcn.Range("C1").ColumnWidth = 29.1
$zhcn.Range("P1").ColumnWidth = 12.8

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusNew
$dede.Range("K2").Value = "2016-08-26 06:47:48"
$dede.Range("P2").Value = ""
$dede.Range("C1").ColumnWidth = 29.1
$dede.Range("P1").ColumnWidth = 12.8
